$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-17 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-18 Wednesday", 2)
$d.Content.Find.Execute("26÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "11÷3=", 2)
$d.Content.Find.Execute("35÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "83÷6=", 2)
$d.Content.Find.Execute("60÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "10÷9=", 2)
$d.Content.Find.Execute("88÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "63÷7=", 2)
$d.Content.Find.Execute("51÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "90÷6=", 2)
$d.Content.Find.Execute("30÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "59÷5=", 2)
$d.Content.Find.Execute("14÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "13÷2=", 2)
$d.Content.Find.Execute("74÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "36÷6=", 2)
$d.Content.Find.Execute("50÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "49÷4=", 2)
$d.Content.Find.Execute("24÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "78÷6=", 2)
$d.Content.Find.Execute("18÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "79÷5=", 2)
$d.Content.Find.Execute("64÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "21÷2=", 2)
$d.Content.Find.Execute("92÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷3=", 2)
$d.Content.Find.Execute("76÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "97÷6=", 2)
$d.Content.Find.Execute("31÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "89÷7=", 2)
$d.Content.Find.Execute("10÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "61÷3=", 2)
$d.Content.Find.Execute("10÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "58÷3=", 2)
$d.Content.Find.Execute("52÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "46÷3=", 2)
$d.Content.Find.Execute("96÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "74÷9=", 2)
$d.Content.Find.Execute("57÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "44÷3=", 2)
$d.Content.Find.Execute("69÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "60÷2=", 2)
$d.Content.Find.Execute("25÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "19÷2=", 2)
$d.Content.Find.Execute("42÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "75÷5=", 2)
$d.Content.Find.Execute("62÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "31÷5=", 2)
$d.Content.Find.Execute("88÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "11÷3=", 2)

Write-Host "Replacements complete"
